# Auto-generated edit script: updates cryptos.xlsx price/volume data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cellRange, [string]$text)
    $origStyle = $cellRange.Style
    $cellRange.NumberFormat = "@"
    $cellRange.Value = $text
    $cellRange.Style = $origStyle
}

Set-TextValue $ws.Range("D2") "27.002.64"
Set-TextValue $ws.Range("E2") "  -0.95%  "
Set-TextValue $ws.Range("D3") "1.619.76"
Set-TextValue $ws.Range("E3") "  -1.53%  "
Set-TextValue $ws.Range("E5") "  -1.33%  "
Set-TextValue $ws.Range("E6") "  -0.13%  "
Set-TextValue $ws.Range("E7") "  -0.05%  "
Set-TextValue $ws.Range("E8") "  -1.17%  "
Set-TextValue $ws.Range("E9") "  -0.42%  "
Set-TextValue $ws.Range("D10") "20.06"
Set-TextValue $ws.Range("E10") "  -0.04%  "
Set-TextValue $ws.Range("E11") "  -0.22%  "
Set-TextValue $ws.Range("D12") "1.846.55"
Set-TextValue $ws.Range("E12") "  -1.51%  "
Set-TextValue $ws.Range("D13") "1.652.11"
Set-TextValue $ws.Range("E13") "  +0.76%  "
Set-TextValue $ws.Range("D14") "4.14"
Set-TextValue $ws.Range("E14") "  -0.19%  "
Set-TextValue $ws.Range("E15") "  -1.16%  "
Set-TextValue $ws.Range("E16") "  -3.65%  "
Set-TextValue $ws.Range("D17") "26.989.44"
Set-TextValue $ws.Range("E17") "  -0.98%  "
Set-TextValue $ws.Range("E18") "  +0.73%  "
Set-TextValue $ws.Range("D19") "213.95"
Set-TextValue $ws.Range("E19") "  -2.87%  "
Set-TextValue $ws.Range("E20") "  -0.07%  "
Set-TextValue $ws.Range("E21") "  -1.93%  "
Set-TextValue $ws.Range("E22") "  -1.28%  "
Set-TextValue $ws.Range("E23") "  -5.63%  "
Set-TextValue $ws.Range("E24") "  -1.34%  "
Set-TextValue $ws.Range("D25") "148.16"
Set-TextValue $ws.Range("E25") "  -0.41%  "
Set-TextValue $ws.Range("E26") "  -0.04%  "
Set-TextValue $ws.Range("E27") "  -1.05%  "
Set-TextValue $ws.Range("E28") "  -1.81%  "
Set-TextValue $ws.Range("D29") "15.56"
Set-TextValue $ws.Range("E29") "  -1.19%  "
Set-TextValue $ws.Range("E30") "  +0.44%  "
Set-TextValue $ws.Range("E31") "  -1.02%  "
Set-TextValue $ws.Range("D32") "3.35"
Set-TextValue $ws.Range("E32") "  -1.51%  "
Set-TextValue $ws.Range("D33") "0.737"
Set-TextValue $ws.Range("E33") "  +33.09%  "
Set-TextValue $ws.Range("E34") "  -0.30%  "
Set-TextValue $ws.Range("D35") "1.342.38"
Set-TextValue $ws.Range("E35") "  +2.61%  "
Set-TextValue $ws.Range("E36") "  -1.03%  "
Set-TextValue $ws.Range("E37") "  -0.52%  "
Set-TextValue $ws.Range("E38") "  +0.38%  "
Set-TextValue $ws.Range("E39") "  -1.68%  "
Set-TextValue $ws.Range("E40") "  -0.07%  "
Set-TextValue $ws.Range("B41") "MXToken"
Set-TextValue $ws.Range("C41") "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue $ws.Range("D41") "2.23"
Set-TextValue $ws.Range("E41") "  +0.03%  "
Set-TextValue $ws.Range("B42") "TrustWalletToken"
Set-TextValue $ws.Range("C42") "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextValue $ws.Range("D42") "0.798"
Set-TextValue $ws.Range("E42") "  -1.71%  "
Set-TextValue $ws.Range("D43") "65.25"
Set-TextValue $ws.Range("E43") "  +5.17%  "
Set-TextValue $ws.Range("D44") "5.32"
Set-TextValue $ws.Range("E44") "  -0.19%  "
Set-TextValue $ws.Range("D45") "1.756.34"
Set-TextValue $ws.Range("E45") "  -1.65%  "
Set-TextValue $ws.Range("D46") "89.70"
Set-TextValue $ws.Range("E46") "  -2.48%  "
Set-TextValue $ws.Range("E47") "  +0.84%  "
Set-TextValue $ws.Range("D48") "0.858"
Set-TextValue $ws.Range("E48") "  +28.35%  "
Set-TextValue $ws.Range("D49") "0.0₆0107"
Set-TextValue $ws.Range("E49") "  -0.58%  "
Set-TextValue $ws.Range("B50") "Algorand"
Set-TextValue $ws.Range("C50") "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue $ws.Range("D50") "0.101"
Set-TextValue $ws.Range("E50") "  +4.61%  "
Set-TextValue $ws.Range("B51") "Cronos"
Set-TextValue $ws.Range("C51") "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue $ws.Range("D51") "0.0516"
Set-TextValue $ws.Range("E51") "  +0.57%  "
